$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: new English translation strings (Part 1 dialogue update) ---
# Values are written in the same order the shared strings were originally
# authored so the resulting sharedStrings.xml table indices line up exactly
# (rows 12 and 28 were filled in last, after the rest, then row 54).
$ws.Range("C4").Value = "Caught it!"
$ws.Range("C6").Value = "There's nothing there anymore…"
$ws.Range("C9").Value = "There's a hole in the wall just the right size for a mouse…"
$ws.Range("C11").Value = "The hole is sealed with ice…"
$ws.Range("C18").Value = "When you're short, it can be a pain to reach high places!`nIt's such a bother to drag a chair everywhere you go!`nHave we got the product for you!"
$ws.Range("C19").Value = "…`nThe rest is torn out."
$ws.Range("C20").Value = "\n<\n[1]>(What was the product… I wonder.)"
$ws.Range("C21").Value = "\n<\n[1]>(A way to reach high-up places…?)"
$ws.Range("C22").Value = "Ornamental swords.`nProbably no good as actual weapons."
$ws.Range("C23").Value = "A piano…"
$ws.Range("C24").Value = "\n<\n[1]>(Oooh… I want to play. But I know I shouldn't.)"
$ws.Range("C25").Value = "\n<\n[1]>(I bet Meria would love to play this.)"
$ws.Range("C26").Value = "\n<\n[1]>(I have to sit beside her and press the keys sometimes.`nShe says her hands are too small to reach a full octave.)"
$ws.Range("C27").Value = "Lily's Diary"
$ws.Range("C29").Value = "Read"
$ws.Range("C30").Value = "Don't read"
$ws.Range("C31").Value = "The clock ticks away…"
$ws.Range("C32").Value = "\n<\n[1]>(I wonder how long I'll be trapped here…)"
$ws.Range("C33").Value = "\n<\n[1]>(I've lost my sense of time…`nI need to get out of here quickly…)"
$ws.Range("C34").Value = "Some flowers are arranged."
$ws.Range("C35").Value = "\n<\n[1]>(What cute flowers. They smell lovely.)"
$ws.Range("C36").Value = "\n<\n[1]>(Anything under here… Nope, no key.)"
$ws.Range("C37").Value = "A houseplant."
$ws.Range("C38").Value = "The bottle is smashed.`nIt smells like alcohol…"
$ws.Range("C39").Value = "\n<\n[1]>(What a waste…)"
$ws.Range("C40").Value = "Lukewarm beer…"
$ws.Range("C41").Value = "\n<\n[1]>(It looks so delicious when my brother drinks it.`nBut I don't actually know how it tastes.)"
$ws.Range("C42").Value = "\n<\n[1]>(I guess they drink it room-temperature here. Disgusting.)"
$ws.Range("C43").Value = "Some wine sits atop the counter…"
$ws.Range("C44").Value = "\n<\n[1]>(I wonder if lots of people used to sit here and drink.)"
$ws.Range("C45").Value = "\n<\n[1]>(A vintage…)"
$ws.Range("C46").Value = "Forks and spoons are all lined up…"
$ws.Range("C50").Value = "How to catch a mouse…"
$ws.Range("C51").Value = "See answer"
$ws.Range("C52").Value = "Don't look"
$ws.Range("C53").Value = "Seal up the mouse hole.`nMethod 1: Push the chest in front`nMethod 2: Seal the hole with ice."
$ws.Range("C12").Value = "\n<\n[1]>(I feel bad about that...)"
$ws.Range("C28").Value = "Lime's Diary`nA Delicious Meal for Semen Slaves"
$ws.Range("C54").Value = "Method 3: Brute force (Approach from above and catch)`nAny of these options is valid."

# --- Wrap text for cells holding multi-line / long translated text ---
$ws.Range("A10").WrapText = $true
$ws.Range("A18").WrapText = $true
$ws.Range("B18").WrapText = $true
$ws.Range("C18").WrapText = $true
$ws.Range("B19").WrapText = $true
$ws.Range("C19").WrapText = $true
$ws.Range("B22").WrapText = $true
$ws.Range("C22").WrapText = $true
$ws.Range("C23").WrapText = $true
$ws.Range("B26").WrapText = $true
$ws.Range("C26").WrapText = $true
$ws.Range("C27").WrapText = $true
$ws.Range("B28").WrapText = $true
$ws.Range("C28").WrapText = $true
$ws.Range("C29").WrapText = $true
$ws.Range("C30").WrapText = $true
$ws.Range("C31").WrapText = $true
$ws.Range("B33").WrapText = $true
$ws.Range("C33").WrapText = $true
$ws.Range("C34").WrapText = $true
$ws.Range("B38").WrapText = $true
$ws.Range("C38").WrapText = $true
$ws.Range("B41").WrapText = $true
$ws.Range("C41").WrapText = $true
$ws.Range("B53").WrapText = $true
$ws.Range("C53").WrapText = $true
$ws.Range("B54").WrapText = $true
$ws.Range("C54").WrapText = $true

# --- Row heights for the wrapped rows ---
$ws.Rows.Item(10).RowHeight = 75
$ws.Rows.Item(18).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(33).RowHeight = 30
$ws.Rows.Item(38).RowHeight = 30
$ws.Rows.Item(41).RowHeight = 30
$ws.Rows.Item(53).RowHeight = 45
$ws.Rows.Item(54).RowHeight = 30

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 76.57
$ws.Columns.Item(2).ColumnWidth = 78.86
$ws.Columns.Item(3).ColumnWidth = 70.71

# --- Sheet view: zoom + final selection, keep gridlines shown ---
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.Zoom = 85
$ws.Range("C52").Select() | Out-Null
